$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column C (shifts old C->D, old D->E)
$ws.Columns("C:C").Insert()

# New header for column C
$ws.Range("C1").Value = "Variance"

# New values in column C: B^2 (typed in C2, then filled down C3:C11)
$ws.Range("C2").Formula = "=B2^2"
$ws.Range("C3:C11").Formula = "=B3^2"

# New average formula for column C (Variance)
$ws.Range("C13").Formula = "=AVERAGE(C2:C11)"

# Remove B14 and C14 (old STD formulas for B and variance columns)
$ws.Range("B14").ClearContents()
$ws.Range("C14").ClearContents()

# Add row 15 with RMS
$ws.Range("A15").Value = "RMS"
$ws.Range("B15").Formula = "=SQRT(C13)"

# Set column C width to match the auto-fit width Excel computed for "Variance"
$ws.Columns("C:C").ColumnWidth = 11.330729166666666

# Set selection to B16
$ws.Range("B16").Select() | Out-Null
